# Actualización automática 2025-11-12 10:30:07
$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M13").Value = 4299.98

# Sheet: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F13").Value = 4299.98
$ws2.Range("F23").Value = 6623.96

# Sheet: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D12").Value = 6260.91
$ws3.Range("E12").Value = 38157.09
$ws3.Range("F12").Value = 0.1409543428339862

$ws3.Range("D14").Value = 6623.96
$ws3.Range("E14").Value = 48775.51101170094
$ws3.Range("F14").Value = 0.1195672066724419
